# spring 24 week 5 inputs
# Append 17 new matchup rows (A1416:D1432) below the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, 8, 7, 12),
    @(7, 7, 5, 13),
    @(2, 17, 4, 3),
    @(5, 8, 4, 12),
    @(1, 14, 3, 6),
    @(6, 5, 5, 15),
    @(3, 5, 5, 15),
    @(8, 14, 5, 6),
    @(3, 16, 4, 4),
    @(2, 14, 4, 6),
    @(3, 4, 2, 16),
    @(6, 12, 9, 8),
    @(2, 16, 4, 4),
    @(2, 7, 1, 13),
    @(3, 12, 4, 8),
    @(5, 16, 6, 4),
    @(4, 4, 2, 16)
)

$startRow = 1416
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Match the author's resulting scroll position / active-cell selection.
$excel.ActiveWindow.ScrollRow = $startRow + 2
$ws.Range("A$($endRow + 1)").Select()
